# Applies the "2025-10-19 01:54 JST" scrape refresh to the "ランサーズ"
# (lancers) sheet: every existing row gets the new capture timestamp,
# two brand-new listings are spliced in (pushing the lower-ranked rows
# down), and the trailing two rows that fall off the bottom of the old
# range are re-created at rows 13-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$timestamp = "2025-10-19 01:54:35"

# row, title, category, price, deadline, url, score, skill-summary (or $null)
$data = @(
    @(2,  "【業務自動化×補助金対応】生成AI活用/日本人モデル画像生成歓迎", "システム開発", "3,000,000 円 ~ 5,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5405834", 395, "🔥AI,Ai ◆自動化"),
    @(3,  "注目 AIプロンプトエンジニア/応答生成トレーナー募集(モバイルアプリ向け)", "システム開発", "5,000 円 ~ 10,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5415842", 330, "🔥AI,Ai ◇アプリ"),
    @(4,  "【急募】ebayAPIを活用したShippingポリシー設定の専門家募集", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5415908", 183, "🔥API"),
    @(5,  "システム開発において活躍できる案件紹介", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5397117", 125, "◆開発,システム開発"),
    @(6,  "イベント出店者管理用ウェブアプリ開発依頼", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5416005", 123, "◆開発 ◇アプリ"),
    @(7,  "仮想通貨取引のBOT作成", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5415610", 118, "★bot"),
    @(8,  "【メンタルヘルス】支援アプリ開発パートナー募集", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5415859", 93, "◆開発 ◇アプリ"),
    @(9,  "【恋愛診断】フルスクラッチ開発・運用サポート募集", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5415986", 75, "◆開発"),
    @(10, "【高額成功報酬】レガシー基幹システムのバイナリ解析とパッチ作成", "システム開発", "1,000,000 円 ~ 3,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5415980", 40, $null),
    @(11, "【動画制作】Sora2での定期生成とTiktok自動投稿依頼", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5415960", 18, $null),
    @(12, "【クリエイティブ】Aurora Creative Lab 外注パートナー募集", "システム開発", "50,000 円 ~ 100,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5415615", 18, $null),
    @(13, "限定公開 限定公開の仕事", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5415804", 13, $null),
    @(14, "【急募】独自ドメインのメール送信エラー解消をお願いいたします", "システム開発", "~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5415841", 10, $null)
)

foreach ($row in $data) {
    $r      = $row[0]
    $title  = $row[1]
    $cat    = $row[2]
    $price  = $row[3]
    $due    = $row[4]
    $url    = $row[5]
    $score  = $row[6]
    $skill  = $row[7]

    $ws.Range("A$r").Value = $timestamp
    $ws.Range("B$r").Value = $title
    $ws.Range("C$r").Value = $cat
    $ws.Range("D$r").Value = $price
    $ws.Range("E$r").Value = $due
    $ws.Range("F$r").Value = $url
    $ws.Range("G$r").Value = $score

    if ($skill -ne $null) {
        $ws.Range("H$r").Value = $skill
    } else {
        $ws.Range("H$r").ClearContents()
    }
}

# Rebuild the hyperlink collection from scratch so the F-column rIds line
# back up with the (now-shifted) URLs instead of leaving stale targets
# behind -- Range-scoped Hyperlinks.Delete() on this engine clears the
# whole-sheet collection anyway, so do it once up front.
$ws.Hyperlinks.Delete()
foreach ($row in $data) {
    $r   = $row[0]
    $url = $row[5]
    [void]$ws.Hyperlinks.Add($ws.Range("F$r"), $url)
}
